# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sorts / refreshes the worker mora records on Hoja1 (rows 16-24):
# period 2503 entries first for all 5 workers, then the 2504 entries for
# the 4 workers that have a second period, then the single 2507 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico

$data = @(
    @{ Row=16; TipoDoc="CC"; Doc="73127722"; Nombre="ALCIDES EDUARDO HARDY MENDOZA";   Periodo="2503"; Mora=7592;  Salario=1423500 },
    @{ Row=17; TipoDoc="CC"; Doc="7920727";  Nombre="JUAN MANUEL ARNEDO HERRERA";      Periodo="2503"; Mora=36062; Salario=1423500 },
    @{ Row=18; TipoDoc="CC"; Doc="73203058"; Nombre="CARLOS RAFAEL PADILLA ESPINOSA"; Periodo="2503"; Mora=7592;  Salario=1423500 },
    @{ Row=19; TipoDoc="CC"; Doc="7252700";  Nombre="JAYR ANTONIO MORENO RESTREPO";   Periodo="2503"; Mora=36062; Salario=1423500 },
    @{ Row=20; TipoDoc="CC"; Doc="79820964"; Nombre="FAVER GIRALDO GARCIA";           Periodo="2503"; Mora=36062; Salario=1423500 },
    @{ Row=21; TipoDoc="CC"; Doc="7920727";  Nombre="JUAN MANUEL ARNEDO HERRERA";      Periodo="2504"; Mora=56940; Salario=1423500 },
    @{ Row=22; TipoDoc="CC"; Doc="7252700";  Nombre="JAYR ANTONIO MORENO RESTREPO";   Periodo="2504"; Mora=56940; Salario=1423500 },
    @{ Row=23; TipoDoc="CC"; Doc="79820964"; Nombre="FAVER GIRALDO GARCIA";           Periodo="2504"; Mora=56940; Salario=1423500 },
    @{ Row=24; TipoDoc="CC"; Doc="79820964"; Nombre="FAVER GIRALDO GARCIA";           Periodo="2507"; Mora=1898;  Salario=1423500 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 2).Value = $rec.TipoDoc
    $ws.Cells.Item($r, 3).Value = $rec.Doc
    $ws.Cells.Item($r, 4).Value = $rec.Nombre
    $ws.Cells.Item($r, 5).Value = $rec.Periodo
    $ws.Cells.Item($r, 6).Value = $rec.Mora
    $ws.Cells.Item($r, 7).Value = $rec.Salario
}
